$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values to reflect the repulled data / mean calc
$ws.Range("F2").Value = -9
$ws.Range("F5").Value = -3
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = -1
$ws.Range("F10").Value = 2
